$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 16:16"

$ws.Range("B7").Value = 27120
$ws.Range("C7").Value = 2913
$ws.Range("E7").Value = 26593

$ws.Range("E13").Value = 4682
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = 243

$ws.Range("B16").Value = 3302
$ws.Range("C16").Value = 310
$ws.Range("E16").Value = 3284

$ws.Range("A21").Value = "Canada"
$ws.Range("B21").Value = 1385
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 1351
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 20

$ws.Range("A22").Value = "Australia"
$ws.Range("B22").Value = 1353
$ws.Range("C22").Value = 281
$ws.Range("D22").Value = 46
$ws.Range("E22").Value = 1300
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 7

$ws.Range("F26").Value = 19

$ws.Range("E49").Value = 181
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 2

$ws.Range("A52").Value = "Hong Kong"
$ws.Range("B52").Value = 317
$ws.Range("C52").Value = 43
$ws.Range("D52").Value = 100
$ws.Range("E52").Value = 213
$ws.Range("F52").Value = 4
$ws.Range("H52").Value = 4

$ws.Range("A53").Value = "Egipto"
$ws.Range("B53").Value = 294
$ws.Range("D53").Value = 42
$ws.Range("E53").Value = 242
$ws.Range("F53").Value = 0
$ws.Range("H53").Value = 10

$ws.Range("A66").Value = "Bulgaria"
$ws.Range("B66").Value = 185
$ws.Range("C66").Value = 22
$ws.Range("D66").Value = 3
$ws.Range("E66").Value = 179
$ws.Range("F66").Value = 3
$ws.Range("H66").Value = 3

$ws.Range("A67").Value = "Eslovaquia"
$ws.Range("B67").Value = 178
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 7
$ws.Range("E67").Value = 171
$ws.Range("F67").Value = 2
$ws.Range("H67").Value = 0

$ws.Range("A68").Value = "San Marino"
$ws.Range("B68").Value = 175
$ws.Range("C68").Value = 15
$ws.Range("D68").Value = 4
$ws.Range("E68").Value = 151
$ws.Range("F68").Value = 13
$ws.Range("H68").Value = 20

$ws.Range("E77").Value = 113
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 1

$ws.Range("A80").Value = "Vietnam"
$ws.Range("B80").Value = 113
$ws.Range("C80").Value = 19
$ws.Range("D80").Value = 17
$ws.Range("E80").Value = 96
$ws.Range("F80").Value = 2
$ws.Range("H80").Value = 0

$ws.Range("A81").Value = "Marruecos"
$ws.Range("B81").Value = 109
$ws.Range("C81").Value = 13
$ws.Range("D81").Value = 3
$ws.Range("E81").Value = 103
$ws.Range("F81").Value = 1
$ws.Range("H81").Value = 3

$ws.Range("A86").Value = "Brunei"
$ws.Range("B86").Value = 88
$ws.Range("C86").Value = 5
$ws.Range("D86").Value = 2
$ws.Range("E86").Value = 86
$ws.Range("F86").Value = 2
$ws.Range("H86").Value = 0

$ws.Range("A87").Value = "Republica de Chipre"
$ws.Range("B87").Value = 84
$ws.Range("D87").Value = 3
$ws.Range("E87").Value = 80
$ws.Range("F87").Value = 3
$ws.Range("H87").Value = 1

$ws.Range("B97").Value = 58
$ws.Range("C97").Value = 4
$ws.Range("E97").Value = 58

$ws.Range("D104").Value = 1
$ws.Range("E104").Value = 46

$ws.Range("A129").Value = "Montenegro"

$ws.Range("A130").Value = "Togo"

$ws.Range("A132").Value = "Kenia"
$ws.Range("B132").Value = 15
$ws.Range("C132").Value = 8
$ws.Range("E132").Value = 15

$ws.Range("A133").Value = "Kirguistan"

$ws.Range("A134").Value = "Barbados"
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 14

$ws.Range("A135").Value = "Costa de Marfil"
$ws.Range("B135").Value = 14
$ws.Range("D135").Value = 1
$ws.Range("E135").Value = 13

$ws.Range("A136").Value = "Maldivas"
$ws.Range("B136").Value = 13
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 3
$ws.Range("E136").Value = 10

$ws.Range("A137").Value = "Tanzania"
$ws.Range("B137").Value = 12
$ws.Range("C137").Value = 6
$ws.Range("E137").Value = 12

$ws.Range("A139").Value = "Etiopia"
$ws.Range("B139").Value = 11
$ws.Range("C139").Value = 2
$ws.Range("E139").Value = 11

$ws.Range("A140").Value = "Mongolia"
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 10

$ws.Range("A141").Value = "Gibraltar"
$ws.Range("B141").Value = 10
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 2
$ws.Range("E141").Value = 8

$ws.Range("A142").Value = "Aruba"
$ws.Range("B142").Value = 8
$ws.Range("C142").Value = 3
$ws.Range("D142").Value = 1
